$wb = $excel.ActiveWorkbook

function Set-StockRows {
    param(
        [string]$SheetName,
        [array]$Rows  # array of arrays: Time, Open, High, Low, Close, Price, Change($), Change(%)
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $startRow = 5
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $startRow + $i
        $data = $Rows[$i]
        $ws.Range("A$r").Value = $data[0]
        $ws.Range("A$r").NumberFormat = "h:mm:ss"
        $ws.Range("B$r").Value = $data[1]
        $ws.Range("C$r").Value = $data[2]
        $ws.Range("D$r").Value = $data[3]
        $ws.Range("E$r").Value = $data[4]
        $ws.Range("F$r").Value = $data[5]
        $ws.Range("G$r").Value = $data[6]
        $ws.Range("H$r").Value = $data[7]
    }
}

# Monster (sheet1)
$monsterRows = @(
    @(0.50209490740740736, 71.760000000000005, 72.819999999999993, 71.22, 71.7, 71.760000000000005, 0.06, 0.08),
    @(0.51635416666666667, 71.760000000000005, 72.819999999999993, 71.22, 71.7, 71.760000000000005, 0.06, 0.08),
    @(0.53057870370370364, 71.760000000000005, 72.819999999999993, 71.22, 71.7, 71.760000000000005, 0.06, 0.08)
)
Set-StockRows "Monster" $monsterRows

# EA (sheet3)
$eaRows = @(
    @(0.50012731481481476, 135.61000000000001, 136.34, 134.51, 134.78, 135.61000000000001, 0.83, 0.62),
    @(0.51503472222222224, 135.61000000000001, 136.34, 134.51, 134.78, 135.61000000000001, 0.83, 0.62),
    @(0.5285185185185185, 135.61000000000001, 136.34, 134.51, 134.78, 135.61000000000001, 0.83, 0.62)
)
Set-StockRows "EA" $eaRows

# Nvidia (sheet5)
$nvidiaRows = @(
    @(0.49677083333333333, 412.8, 409, 398.62, 394.87, 412.8, 4.16, 1.02),
    @(0.51265046296296302, 413.52, 409, 398.62, 394.87, 413.52, 4.88, 1.19),
    @(0.52491898148148153, 414.89, 409, 398.62, 394.87, 414.89, 6.25, 1.53)
)
Set-StockRows "Nvidia" $nvidiaRows

# Microsoft (sheet6)
$msftRows = @(
    @(0.49525462962962963, 214.2, 213.26, 208.69, 208.25, 214.2, 1.37, 0.64),
    @(0.51143518518518516, 214.54, 213.26, 208.69, 208.25, 214.54, 1.71, 0.8),
    @(0.52339120370370373, 215.14, 213.26, 208.69, 208.25, 215.14, 2.31, 1.0900000000000001)
)
Set-StockRows "Microsoft" $msftRows
